$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "arab_title_script"
$ws.Range("C1").Value = "arab_title"
$ws.Range("D1").Value = "chinese_title"
$ws.Range("E1").Value = "author"
$ws.Range("F1").Value = "assembler"
$ws.Range("G1").Value = "editor"
$ws.Range("H1").Value = "scrivener"
$ws.Range("I1").Value = "translator"
$ws.Range("J1").Value = "type"
$ws.Range("K1").Value = "place"
$ws.Range("L1").Value = "publisher"
$ws.Range("M1").Value = "year"
$ws.Range("N1").Value = "stand_year"
$ws.Range("O1").Value = "language"
$ws.Range("P1").Value = "num_pages"
$ws.Range("Q1").Value = "description"
$ws.Range("R1").Value = "notes"
